$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = "carne bovina moída"
$ws.Cells.Item(7, 2).Value = "comida"
$ws.Cells.Item(7, 3).Value = "carne"
$ws.Cells.Item(7, 4).Value = 250.0
$ws.Cells.Item(7, 5).Value = 20.0
$ws.Cells.Item(7, 6).Value = 17.0
$ws.Cells.Item(7, 7).Value = 0.0

$ws.Cells.Item(8, 1).Value = "macarrão"
$ws.Cells.Item(8, 2).Value = "comida"
$ws.Cells.Item(8, 3).Value = "carboidrato"
$ws.Cells.Item(8, 4).Value = 131.0
$ws.Cells.Item(8, 5).Value = 1.1
$ws.Cells.Item(8, 6).Value = 5.0
$ws.Cells.Item(8, 7).Value = 25.0

$ws.Cells.Item(9, 1).Value = "chocolate 70%"
$ws.Cells.Item(9, 2).Value = "comida"
$ws.Cells.Item(9, 3).Value = "doce"
$ws.Cells.Item(9, 4).Value = 598.0
$ws.Cells.Item(9, 5).Value = 42.6
$ws.Cells.Item(9, 6).Value = 7.6
$ws.Cells.Item(9, 7).Value = 45.9

$ws.Cells.Item(10, 1).Value = "chocolate 70%"
$ws.Cells.Item(10, 2).Value = "comida"
$ws.Cells.Item(10, 3).Value = "doce"
$ws.Cells.Item(10, 4).Value = 535.0
$ws.Cells.Item(10, 5).Value = 29.7
$ws.Cells.Item(10, 6).Value = 7.6
$ws.Cells.Item(10, 7).Value = 59.4

$ws.Cells.Item(11, 1).Value = "amêndoa"
$ws.Cells.Item(11, 2).Value = "comida"
$ws.Cells.Item(11, 3).Value = "carboidrato"
$ws.Cells.Item(11, 4).Value = 576.0
$ws.Cells.Item(11, 5).Value = 49.4
$ws.Cells.Item(11, 6).Value = 21.2
$ws.Cells.Item(11, 7).Value = 21.6

$ws.Cells.Item(12, 1).Value = "castanha de caju"
$ws.Cells.Item(12, 2).Value = "comida"
$ws.Cells.Item(12, 3).Value = "carboidrato"
$ws.Cells.Item(12, 4).Value = 570.0
$ws.Cells.Item(12, 5).Value = 46.3
$ws.Cells.Item(12, 6).Value = 18.5
$ws.Cells.Item(12, 7).Value = 29.1

$ws.Cells.Item(13, 1).Value = "nozes"
$ws.Cells.Item(13, 2).Value = "comida"
$ws.Cells.Item(13, 3).Value = "carboidrato"
$ws.Cells.Item(13, 4).Value = 654.0
$ws.Cells.Item(13, 5).Value = 65.2
$ws.Cells.Item(13, 6).Value = 15.2
$ws.Cells.Item(13, 7).Value = 13.7

$ws.Cells.Item(14, 1).Value = "capuccino"
$ws.Cells.Item(14, 2).Value = "bebida"
$ws.Cells.Item(14, 3).Value = "-"
$ws.Cells.Item(14, 4).Value = 30.0
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 4.0

$ws.Cells.Item(15, 1).Value = "pão branco"
$ws.Cells.Item(15, 2).Value = "comida"
$ws.Cells.Item(15, 3).Value = "carboidrato"
$ws.Cells.Item(15, 4).Value = 256.0
$ws.Cells.Item(15, 5).Value = 3.2
$ws.Cells.Item(15, 6).Value = 9.0
$ws.Cells.Item(15, 7).Value = 49.0

$ws.Cells.Item(16, 1).Value = "pão integral"
$ws.Cells.Item(16, 2).Value = "comida"
$ws.Cells.Item(16, 3).Value = "carboidrato"
$ws.Cells.Item(16, 4).Value = 247.0
$ws.Cells.Item(16, 5).Value = 4.2
$ws.Cells.Item(16, 6).Value = 9.6
$ws.Cells.Item(16, 7).Value = 41.4

$ws.Cells.Item(17, 1).Value = "queijo mussarela"
$ws.Cells.Item(17, 2).Value = "comida"
$ws.Cells.Item(17, 3).Value = "proteína"
$ws.Cells.Item(17, 4).Value = 280.0
$ws.Cells.Item(17, 5).Value = 17.0
$ws.Cells.Item(17, 6).Value = 27.0
$ws.Cells.Item(17, 7).Value = 3.1

# Apply the same cell formatting (style) used by the existing data rows
$ws.Range("A2:G2").Copy()
$ws.Range("A7:G17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
